$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force columns D and E to text format while writing, so values like "353.74" or
# "52.210.30" are stored as strings (matching the source data) rather than being
# auto-parsed into numbers by Excel. ClearFormats afterwards restores the original
# (unset/General) cell style so only the cell VALUES change, not their formatting.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$ws.Range('D2').Value = '52.210.30'
$ws.Range('E2').Value = '  +0.80%  '
$ws.Range('D3').Value = '2.907.16'
$ws.Range('E3').Value = '  +3.69%  '
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').Value = '353.74'
$ws.Range('E5').Value = '  +0.05%  '
$ws.Range('D6').Value = '114.05'
$ws.Range('E6').Value = '  +1.60%  '
$ws.Range('D7').Value = '0.558'
$ws.Range('E7').Value = '  +0.12%  '
$ws.Range('E8').Value = '  +0.08%  '
$ws.Range('E9').Value = '  -0.11%  '
$ws.Range('D10').Value = '39.73'
$ws.Range('E10').Value = '  -1.16%  '
$ws.Range('D11').Value = '0.0868'
$ws.Range('E11').Value = '  +3.51%  '
$ws.Range('E12').Value = '  +0.80%  '
$ws.Range('D13').Value = '19.81'
$ws.Range('E13').Value = '  -0.44%  '
$ws.Range('D14').Value = '7.74'
$ws.Range('E14').Value = '  -0.37%  '
$ws.Range('D15').Value = '3.372.56'
$ws.Range('E15').Value = '  +3.95%  '
$ws.Range('D16').Value = '2.907.79'
$ws.Range('E16').Value = '  +3.65%  '
$ws.Range('D17').Value = '0.987'
$ws.Range('E17').Value = '  +3.97%  '
$ws.Range('D18').Value = '52.294.07'
$ws.Range('E18').Value = '  +0.91%  '
$ws.Range('E19').Value = '  +2.05%  '
$ws.Range('D20').Value = '7.62'
$ws.Range('E20').Value = '  +0.09%  '
$ws.Range('D21').Value = '14.10'
$ws.Range('E21').Value = '  +4.27%  '
$ws.Range('D22').Value = '0.0₃0978'
$ws.Range('E22').Value = '  +0.45%  '
$ws.Range('D23').Value = '71.09'
$ws.Range('E23').Value = '  +1.08%  '
$ws.Range('D24').Value = '269.28'
$ws.Range('E24').Value = '  +0.76%  '
$ws.Range('D25').Value = '2.82'
$ws.Range('E25').Value = '  +1.98%  '
$ws.Range('D26').Value = '0.182'
$ws.Range('E26').Value = '  +12.62%  '
$ws.Range('D27').Value = '26.77'
$ws.Range('E27').Value = '  +2.38%  '
$ws.Range('E28').Value = '  -0.20%  '
$ws.Range('D29').Value = '10.65'
$ws.Range('E29').Value = '  +2.53%  '
$ws.Range('D30').Value = '0.103'
$ws.Range('E30').Value = '  +15.03%  '
$ws.Range('D31').Value = '6.79'
$ws.Range('E31').Value = '  +11.05%  '
$ws.Range('D32').Value = '37.56'
$ws.Range('E32').Value = '  -4.30%  '
$ws.Range('E33').Value = '  -0.51%  '
$ws.Range('D34').Value = '6.10'
$ws.Range('E34').Value = '  +10.46%  '
$ws.Range('D35').Value = '53.06'
$ws.Range('E35').Value = '  +1.65%  '
$ws.Range('D36').Value = '0.0452'
$ws.Range('E36').Value = '  +0.51%  '
$ws.Range('E37').Value = '  +0.00%  '
$ws.Range('D38').Value = '3.32'
$ws.Range('E38').Value = '  +4.51%  '
$ws.Range('D39').Value = '18.82'
$ws.Range('E39').Value = '  -0.78%  '
$ws.Range('D40').Value = '2.04'
$ws.Range('E40').Value = '  +1.69%  '
$ws.Range('D41').Value = '2.74'
$ws.Range('E41').Value = '  +8.79%  '
$ws.Range('E42').Value = '  +1.36%  '
$ws.Range('D43').Value = '23.11'
$ws.Range('E43').Value = '  +5.41%  '
$ws.Range('D44').Value = '118.06'
$ws.Range('E44').Value = '  -1.38%  '
$ws.Range('D45').Value = '2.17'
$ws.Range('E45').Value = '  -2.42%  '
$ws.Range('E46').Value = '  +2.08%  '
$ws.Range('D47').Value = '3.53'
$ws.Range('E47').Value = '  -0.10%  '
$ws.Range('D48').Value = '2.179.83'
$ws.Range('E48').Value = '  +3.15%  '
$ws.Range('E49').Value = '  +17.69%  '
$ws.Range('E50').Value = '  +12.23%  '
$ws.Range('D51').Value = '0.954'
$ws.Range('E51').Value = '  -3.13%  '

$dataRange.ClearFormats()
